$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 182
$ws.Range("D2").Value = 95.51

$ws.Range("C3").Value = 184

$ws.Range("C4").Value = 177
$ws.Range("D4").Value = 104.14

$ws.Range("C5").Value = 236
$ws.Range("D5").Value = 102.2

$ws.Range("C6").Value = 201
$ws.Range("D6").Value = 99.94

$ws.Range("C7").Value = 87
$ws.Range("D7").Value = 42.58

$ws.Range("D8").Value = 35
